$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-13 00:45:02"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-13 00:44:52"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/1bd89b1ed7b06cfcf3072f5a0b661f821ee8d820/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/e985fbc7801c19f744490ae443fc2551613ce127/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = 39.17

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "'False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-13 00:45:02"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/1bd89b1ed7b06cfcf3072f5a0b661f821ee8d820/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/e985fbc7801c19f744490ae443fc2551613ce127/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 39.17
